# Append: 2026-01-21 06:33 JST
# - Refresh the "取得日時" timestamp on the five surviving rows
# - Replace the two newest listings (rows 5 & 6) with the latest scrape results
# - Drop the now-stale listings that previously occupied rows 7-16
# - Narrow column D slightly and fix up the hyperlink relationships

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$newTimestamp = "2026-01-21 06:33:53"

# Remove the old rows (7-16) first, while row numbers still match the
# pre-edit layout. Everything below row 6 disappears.
$ws.Range("A7:H16").EntireRow.Delete()

# Refresh the capture timestamp for every remaining data row (2-6).
$ws.Range("A2").Value = $newTimestamp
$ws.Range("A3").Value = $newTimestamp
$ws.Range("A4").Value = $newTimestamp
$ws.Range("A5").Value = $newTimestamp
$ws.Range("A6").Value = $newTimestamp

# Row 5 becomes a brand-new listing.
$ws.Range("B5").Value = "【急募】野球スコアボードシステム開発のフリーランス募集"
$ws.Range("D5").Value = "50,000 円 ~ 100,000 円 / 固定"
$ws.Range("F5").Value = "https://www.lancers.jp/work/detail/5475665"
$ws.Range("G5").Value = 118
$ws.Range("H5").Value = "◆開発,システム開発"

# Row 6 becomes a brand-new listing too; it no longer carries an H value.
$ws.Range("B6").Value = "【急募】CSVデータをワードに自動入力するスキルをお持ちの方"
$ws.Range("D6").Value = "10,000 円 ~ 20,000 円 / 固定"
$ws.Range("F6").Value = "https://www.lancers.jp/work/detail/5475924"
$ws.Range("G6").Value = 10
$ws.Range("H6").ClearContents()

# Column D gets a touch narrower.
$ws.Columns.Item(4).ColumnWidth = 27.17

# The row deletion above leaves the worksheet's hyperlink collection
# pointing at relationship ids for rows that no longer exist, and F5/F6
# now need new target URLs. Rebuild the hyperlinks for the five surviving
# URL cells from scratch so the relationships line up with rId1..rId5.
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("F2"), "https://www.lancers.jp/work/detail/5455098")
$ws.Hyperlinks.Add($ws.Range("F3"), "https://www.lancers.jp/work/detail/5445159")
$ws.Hyperlinks.Add($ws.Range("F4"), "https://www.lancers.jp/work/detail/5445154")
$ws.Hyperlinks.Add($ws.Range("F5"), "https://www.lancers.jp/work/detail/5475665")
$ws.Hyperlinks.Add($ws.Range("F6"), "https://www.lancers.jp/work/detail/5475924")

# Hyperlinks.Add() re-styles each cell with a freshly minted (but
# equivalent) "Hyperlink" style entry. Re-apply the named style explicitly
# so these cells keep using the workbook's original Hyperlink cellXf
# instead of accumulating duplicate style records.
$ws.Range("F2").Style = "Hyperlink"
$ws.Range("F3").Style = "Hyperlink"
$ws.Range("F4").Style = "Hyperlink"
$ws.Range("F5").Style = "Hyperlink"
$ws.Range("F6").Style = "Hyperlink"
